$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values that changed in row 3
$ws.Range("E3").Value = 5
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 13

# Reflect the saved selection/active cell (E3) as recorded in the sheet view
$ws.Range("E3").Select() | Out-Null
